$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "CELL NAMES" column (C2:C29) from short "Cxx" abbreviations
# to the fully spelled out "Cell xx" form.
$cellNames = @(
    "STA-1E Cell 4N",
    "STA-1E Cell 4S",
    "STA-1E Cell 6",
    "STA-1W Cell 1A",
    "STA-1W Cell 1B",
    "STA-1W Cell 2A",
    "STA-1W Cell 2B",
    "STA-1W Cell 3",
    "STA-1W Cell 4",
    "STA-1W Cell 5A",
    "STA-1W Cell 5B",
    "STA-2 Cell 2",
    "STA-2 Cell 3",
    "STA-2 Cell 4",
    "STA-2 Cell 5",
    "STA-2 Cell 6",
    "STA-2 Cell 8",
    "STA-34 Cell 1B",
    "STA-34 Cell 2B",
    "STA-34 Cell 3B",
    "STA-56 Cell 1A",
    "STA-56 Cell 1B",
    "STA-56 Cell 2A",
    "STA-56 Cell 2B",
    "STA-56 Cell 3A",
    "STA-56 Cell 3B",
    "STA-56 Cell 4B",
    "STA-56 Cell 5B"
)

for ($i = 0; $i -lt $cellNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cellNames[$i]
}

# Widen column C (the newly lengthened "CELL NAMES" text) so the longer
# labels fit / are centered nicely, and move the selection down to the
# last cell of the column (C29), matching the saved UI state.
$ws.Columns.Item(3).ColumnWidth = 14.3

[void]$ws.Range("C29").Select()
